$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# Row 5 (Treatment/Campaign/MDA, age 50-65) is being removed: its upper age
# bound (65) is merged into row 4 (Treatment/Campaign/MDA, age 15-50), which
# becomes age 15-65. Row 5 is then deleted entirely and all following rows
# shift up by one.
$ws.Range("G4").Value = 65
$ws.Rows.Item(5).Delete()

# Update the active selection to reflect where Excel would leave the cursor
# after this edit (matches the saved view state in the target file).
$ws.Activate()
$ws.Range("H13").Select()
